$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column)
$ws.Columns("N:N").Insert()

# The new column should take on the same width as column M (matches Excel's
# default behavior of copying the left neighbor's formatting on insert)
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Update the active selection to reflect the new layout (P6 selected, per diff)
$ws.Range("P6").Select()
